# Push choices sheet display.text into display.title.text
# Also push display.image into display.title.image
$wb = $excel.ActiveWorkbook

# --- choices sheet: rename the display.text column header to display.title.text ---
$choices = $wb.Worksheets.Item("choices")
$choices.Range("C1").Value = "display.title.text"
$choices.Columns.Item(3).ColumnWidth = 17.3
$choices.Range("C2").Select()

# --- queries sheet: normalize the linked-table selection / selectionArgs ---
$queries = $wb.Worksheets.Item("queries")
$queries.Range("E2").Value = "1=1"
$queries.Range("F2").Value = "[]"
$queries.Range("E3").Value = "1=1"
$queries.Range("F3").Value = "[]"
$queries.Range("F4").Select()

# restore queries as the active/selected tab (matches original tabSelected state)
$queries.Activate()
